$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) - add E1 and F1, matching the style of existing header cell D1
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4

$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2 - update existing B2:D2 values, add E2 and F2
$ws.Range("B2").Value = -0.01014579726076258
$ws.Range("C2").Value = 0.02666081550071997
$ws.Range("D2").Value = -0.0767144750377244
$ws.Range("E2").Value = -0.03678720286518841
$ws.Range("F2").Value = 0.08551102827500009

# Row 3 - update existing B3:D3 values, add E3 and F3
$ws.Range("B3").Value = 0.002350303899347969
$ws.Range("C3").Value = 0.05708249221189834
$ws.Range("D3").Value = 0.1002612127618704
$ws.Range("E3").Value = -0.06442888459132723
$ws.Range("F3").Value = -0.09187687035283353
